$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Связь"
$ws.Range("B1").Value = "Одежда/обувь"
$ws.Range("C1").Value = "Копилка"
$ws.Range("D1").Value = "Гигиена"
$ws.Range("E1").Value = "На всякий случай"
$ws.Range("F1").Value = "Еда"

$ws.Range("A2").Value = 140
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 3500

$ws.Range("F10").Select()
